$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for existing rows 2-28
# from 45450 (2024-06-07) to 45451 (2024-06-08).
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45451
}

# Row 28 gains an explicit row height (15, custom) in the target file.
$ws.Rows.Item(28).RowHeight = 15

# Append a new data row (29) with the new case.
$ws.Cells.Item(29, 1).Value = "A 23041-2024"

$ws.Cells.Item(29, 2).Value = 45450
$ws.Cells.Item(29, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(29, 3).Value = 45451
$ws.Cells.Item(29, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"

$ws.Cells.Item(29, 7).Value = 5.9
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0

$ws.Cells.Item(29, 18).Value = "'"
$ws.Cells.Item(29, 18).ClearFormats()
$ws.Cells.Item(29, 18).WrapText = $true
